# Sync attendance_reports: normalize the "Recorded By" (column G) token order.
# Two exact legacy strings are re-ordered workbook-wide:
#   "dnasr281@gmail.com, System"            -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com"   -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
}
